$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer / "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Octubre de 2020 a las 18:53"

# --- Reorder "Montserrat" / "Islas Malvinas" rows (their statistics swap places) ---
$ws.Range("A216").Value = "Montserrat"
$ws.Range("A217").Value = "Islas Malvinas"

# --- Updated country statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 8770951
$ws.Range("C4").Value = 23998
$ws.Range("D4").Value = 5705323
$ws.Range("E4").Value = 2836058
$ws.Range("G4").Value = 286
$ws.Range("H4").Value = 229570

# India (row 5)
$ws.Range("B5").Value = 7845328
$ws.Range("C5").Value = 31660
$ws.Range("D5").Value = 7054252
$ws.Range("E5").Value = 672756
$ws.Range("G5").Value = 328
$ws.Range("H5").Value = 118320

# Brasil (row 6)
$ws.Range("B6").Value = 5358498
$ws.Range("C6").Value = 2848
$ws.Range("E6").Value = 404061
$ws.Range("G6").Value = 37
$ws.Range("H6").Value = 156565

# Reino Unido (row 14)
$ws.Range("B14").Value = 854010
$ws.Range("C14").Value = 23012
$ws.Range("G14").Value = 174
$ws.Range("H14").Value = 44745

# Alemania (row 20)
$ws.Range("B20").Value = 423067
$ws.Range("C20").Value = 5717
$ws.Range("E20").Value = 98870
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 10097

# Turquia (row 24)
$ws.Range("B24").Value = 359784
$ws.Range("C24").Value = 2091
$ws.Range("D24").Value = 313093
$ws.Range("E24").Value = 36964
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = 9727

# Israel (row 28)
$ws.Range("B28").Value = 309374
$ws.Range("C28").Value = 534
$ws.Range("D28").Value = 291130
$ws.Range("E28").Value = 15878
$ws.Range("G28").Value = 37
$ws.Range("H28").Value = 2366

# Canada (row 33)
$ws.Range("B33").Value = 213727
$ws.Range("C33").Value = 1995
$ws.Range("D33").Value = 179537
$ws.Range("E33").Value = 24270
$ws.Range("G33").Value = 32
$ws.Range("H33").Value = 9920

# Ecuador (row 36)
$ws.Range("B36").Value = 159614
$ws.Range("C36").Value = 1344
$ws.Range("E36").Value = 12885
$ws.Range("G36").Value = 14
$ws.Range("H36").Value = 12542

# Guatemala (row 49)
$ws.Range("B49").Value = 104632
$ws.Range("C49").Value = 730
$ws.Range("D49").Value = 93880
$ws.Range("E49").Value = 7143
$ws.Range("G49").Value = 15
$ws.Range("H49").Value = 3609

# Kenia (row 76)
$ws.Range("B76").Value = 48790
$ws.Range("C76").Value = 947
$ws.Range("D76").Value = 33876
$ws.Range("E76").Value = 14018
$ws.Range("G76").Value = 12
$ws.Range("H76").Value = 896

# Grecia (row 87)
$ws.Range("B87").Value = 29992
$ws.Range("C87").Value = 935
$ws.Range("E87").Value = 19439
$ws.Range("G87").Value = 5
$ws.Range("H87").Value = 564

# Republica de Macedonia (row 90)
$ws.Range("B90").Value = 26394
$ws.Range("C90").Value = 403
$ws.Range("D90").Value = 18430
$ws.Range("E90").Value = 7067
$ws.Range("G90").Value = 14
$ws.Range("H90").Value = 897

# Montserrat (row 216, formerly Islas Malvinas) / Islas Malvinas (row 217, formerly Montserrat)
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0
